$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("books")

# Header: "ISBNs (kinda)" -> "ISBNs"
$ws.Range("C1").Value = "ISBNs"

# Rows 11/12: replace placeholder text in column C with the real ISBN number,
# matching the numeric style already used by the rest of column C (copy format from C10).
$ws.Range("C10").Copy()
$ws.Range("C11:C12").PasteSpecial(-4122)
$ws.Range("C11").Value = 1285159454
$ws.Range("C12").Value = 1285159454
$excel.CutCopyMode = $false

# Row 17: remove the stray "baking pi" / 3.14 / rating entries, keep only the ISBN in column C.
$ws.Range("A17").Clear()
$ws.Range("B17").Clear()
$ws.Range("D17").Clear()
$ws.Range("C17").Value = 1285159454

# Update the sheet's view/selection.
[void]$ws.Range("C1").Select()
